# "updating prices at 11:31:01"
# Append a new data row (row 25) to the driver prices sheet. The new
# row carries forward the same price values as the previous row (24),
# only the timestamp in column A is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prevRow = 24
$newRow = 25

# Copy the whole previous row (values + formatting, including the
# "empty" cell in column N) down into the new row.
$src = $ws.Range("A$prevRow`:V$prevRow")
$dst = $ws.Range("A$newRow`:V$newRow")
$src.Copy($dst)

# Stamp the new row with its own timestamp.
$ws.Cells.Item($newRow, 1).Value = 44042.95833333334
